$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlCenter = -4108
$xlLeft = -4131

$ws.Range("A26").Value = 25
$ws.Range("A26").HorizontalAlignment = $xlCenter

$ws.Range("B26").Value = "Bottom View of a Binary Tree"
$ws.Range("B26").HorizontalAlignment = $xlLeft

$ws.Range("D26").Value = "Tree"
$ws.Range("D26").HorizontalAlignment = $xlCenter

$ws.Range("E26").Value = "medium"
$ws.Range("E26").HorizontalAlignment = $xlCenter

$ws.Range("F26").Value = "GeeksForGeeks"
$ws.Range("F26").HorizontalAlignment = $xlCenter
